$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated sval data (filtered save games), rows 2-12, columns B-G
$data = @{
    2  = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 0, 8.656069925401464)
    3  = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 0, 6.82939032824165)
    4  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248)
    5  = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 24.14949828602258)
    6  = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 0, 8.656069925401464)
    7  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 6.15379541431027)
    8  = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 0, 6.82939032824165)
    9  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
    10 = @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1, 2.213936997104367)
    11 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 3.755628166162433)
    12 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
